$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '27.079.84'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E2').Value = '  +0.36%  '
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.825.93'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E3').Value = '  +0.35%  '
$ws.Range('E4').Value = '  +0.67%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '312.54'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +0.74%  '
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('E7').Value = '  +0.33%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3649'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -0.44%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.07381'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  +0.41%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.8795'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  +0.72%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '20.24'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -0.20%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '1.889.41'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  +4.59%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.07325'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  +3.01%  '
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '93.25'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E14').Value = '  +1.73%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '5.362'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -0.89%  '
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '6.524'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  +0.21%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '1.007'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.29%  '
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.000008716'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -0.17%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '1.009'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.68%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '27.564.37'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  +2.08%  '
$ws.Range('E21').Value = '  -0.27%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '5.238'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -0.95%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '10.62'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +0.17%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.093.00'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  +2.86%  '
$ws.Range('E25').Value = '  -0.57%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '151.58'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('E27').Value = '  +0.39%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.126'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -1.15%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '5.174'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -1.48%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '116.39'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  -0.62%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.08936'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('E32').Value = '  +0.08%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.7418'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -2.24%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '4.510'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -0.08%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.947'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  +0.77%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '1.009'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  +0.56%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '1.088'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  -1.32%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '0.05292'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -0.31%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01941'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -0.22%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.963'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -0.31%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.412'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +2.00%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '7.237'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +0.76%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.5243'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -0.92%  '
$ws.Range('E44').Value = '  -0.84%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '8.384'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -0.95%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '0.4870'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  -0.30%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '10.39'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -0.60%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '1.009'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  +0.63%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '104.29'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('E50').Value = '  -0.90%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.06276'
$c.NumberFormat = 'General'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  -0.38%  '
